$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 08:47:19"
$ws1.Range("A3").Value = "Total filas: 128"

$ws1.Cells.Item(92, 1).Value = "08:47:19"
$ws1.Cells.Item(92, 2).Value = "08:54"
$ws1.Cells.Item(92, 3).Value = "10_OLMOS"
$ws1.Cells.Item(92, 4).Value = 7
$ws1.Cells.Item(92, 5).Value = "LP1912"

$ws1.Cells.Item(93, 1).Value = "08:16:48"
$ws1.Cells.Item(93, 2).Value = "08:55"
$ws1.Cells.Item(93, 3).Value = "10_OLMOS"
$ws1.Cells.Item(93, 4).Value = 39
$ws1.Cells.Item(93, 5).Value = "LP1912"

$ws1.Cells.Item(94, 1).Value = "07:14:27"
$ws1.Cells.Item(94, 2).Value = "09:01"
$ws1.Cells.Item(94, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(94, 4).Value = 107
$ws1.Cells.Item(94, 5).Value = "LP1912"

$ws1.Cells.Item(95, 1).Value = "07:44:08"
$ws1.Cells.Item(95, 2).Value = "09:02"
$ws1.Cells.Item(95, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(95, 4).Value = 78
$ws1.Cells.Item(95, 5).Value = "LP1912"

$ws1.Cells.Item(96, 1).Value = "07:57:27"
$ws1.Cells.Item(96, 2).Value = "09:03"
$ws1.Cells.Item(96, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(96, 4).Value = 66
$ws1.Cells.Item(96, 5).Value = "LP1912"

$ws1.Cells.Item(97, 1).Value = "08:16:48"
$ws1.Cells.Item(97, 2).Value = "09:04"
$ws1.Cells.Item(97, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(97, 4).Value = 48
$ws1.Cells.Item(97, 5).Value = "LP1912"

$ws1.Cells.Item(98, 1).Value = "08:47:19"
$ws1.Cells.Item(98, 2).Value = "09:05"
$ws1.Cells.Item(98, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(98, 4).Value = 18
$ws1.Cells.Item(98, 5).Value = "LP1912"

$ws1.Cells.Item(99, 1).Value = "07:14:27"
$ws1.Cells.Item(99, 2).Value = "09:07"
$ws1.Cells.Item(99, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(99, 4).Value = 113
$ws1.Cells.Item(99, 5).Value = "LP1912"

$ws1.Cells.Item(100, 1).Value = "07:44:08"
$ws1.Cells.Item(100, 2).Value = "09:08"
$ws1.Cells.Item(100, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(100, 4).Value = 84
$ws1.Cells.Item(100, 5).Value = "LP1912"

$ws1.Cells.Item(101, 1).Value = "07:57:27"
$ws1.Cells.Item(101, 2).Value = "09:09"
$ws1.Cells.Item(101, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(101, 4).Value = 72
$ws1.Cells.Item(101, 5).Value = "LP1912"

$ws1.Cells.Item(102, 1).Value = "07:14:27"
$ws1.Cells.Item(102, 2).Value = "09:10"
$ws1.Cells.Item(102, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(102, 4).Value = 116
$ws1.Cells.Item(102, 5).Value = "LP1912"

$ws1.Cells.Item(103, 1).Value = "07:44:08"
$ws1.Cells.Item(103, 2).Value = "09:11"
$ws1.Cells.Item(103, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(103, 4).Value = 87
$ws1.Cells.Item(103, 5).Value = "LP1912"

$ws1.Cells.Item(104, 1).Value = "08:33:47"
$ws1.Cells.Item(104, 2).Value = "09:13"
$ws1.Cells.Item(104, 3).Value = "10_OLMOS"
$ws1.Cells.Item(104, 4).Value = 40
$ws1.Cells.Item(104, 5).Value = "LP1912"

$ws1.Cells.Item(105, 1).Value = "07:44:08"
$ws1.Cells.Item(105, 2).Value = "09:14"
$ws1.Cells.Item(105, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(105, 4).Value = 90
$ws1.Cells.Item(105, 5).Value = "LP1912"

$ws1.Cells.Item(106, 1).Value = "07:57:27"
$ws1.Cells.Item(106, 2).Value = "09:16"
$ws1.Cells.Item(106, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(106, 4).Value = 79
$ws1.Cells.Item(106, 5).Value = "LP1912"

$ws1.Cells.Item(107, 1).Value = "07:44:08"
$ws1.Cells.Item(107, 2).Value = "09:17"
$ws1.Cells.Item(107, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(107, 4).Value = 93
$ws1.Cells.Item(107, 5).Value = "LP1912"

$ws1.Cells.Item(108, 1).Value = "07:44:08"
$ws1.Cells.Item(108, 2).Value = "09:21"
$ws1.Cells.Item(108, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(108, 4).Value = 97
$ws1.Cells.Item(108, 5).Value = "LP1912"

$ws1.Cells.Item(109, 1).Value = "07:57:27"
$ws1.Cells.Item(109, 2).Value = "09:22"
$ws1.Cells.Item(109, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(109, 4).Value = 85
$ws1.Cells.Item(109, 5).Value = "LP1912"

$ws1.Cells.Item(110, 1).Value = "07:57:27"
$ws1.Cells.Item(110, 2).Value = "09:22"
$ws1.Cells.Item(110, 3).Value = "17_ROMERO"
$ws1.Cells.Item(110, 4).Value = 85
$ws1.Cells.Item(110, 5).Value = "LP1912"

$ws1.Cells.Item(111, 1).Value = "07:44:08"
$ws1.Cells.Item(111, 2).Value = "09:23"
$ws1.Cells.Item(111, 3).Value = "17_ROMERO"
$ws1.Cells.Item(111, 4).Value = 99
$ws1.Cells.Item(111, 5).Value = "LP1912"

$ws1.Cells.Item(112, 1).Value = "07:57:27"
$ws1.Cells.Item(112, 2).Value = "09:23"
$ws1.Cells.Item(112, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(112, 4).Value = 86
$ws1.Cells.Item(112, 5).Value = "LP1912"

$ws1.Cells.Item(113, 1).Value = "08:47:19"
$ws1.Cells.Item(113, 2).Value = "09:23"
$ws1.Cells.Item(113, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(113, 4).Value = 36
$ws1.Cells.Item(113, 5).Value = "LP1912"

$ws1.Cells.Item(114, 1).Value = "07:44:08"
$ws1.Cells.Item(114, 2).Value = "09:24"
$ws1.Cells.Item(114, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(114, 4).Value = 100
$ws1.Cells.Item(114, 5).Value = "LP1912"

$ws1.Cells.Item(115, 1).Value = "08:16:48"
$ws1.Cells.Item(115, 2).Value = "09:29"
$ws1.Cells.Item(115, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(115, 4).Value = 73
$ws1.Cells.Item(115, 5).Value = "LP1912"

$ws1.Cells.Item(116, 1).Value = "07:44:08"
$ws1.Cells.Item(116, 2).Value = "09:32"
$ws1.Cells.Item(116, 3).Value = "15_ABASTO"
$ws1.Cells.Item(116, 4).Value = 108
$ws1.Cells.Item(116, 5).Value = "LP1912"

$ws1.Cells.Item(117, 1).Value = "07:44:08"
$ws1.Cells.Item(117, 2).Value = "09:33"
$ws1.Cells.Item(117, 3).Value = "10_OLMOS"
$ws1.Cells.Item(117, 4).Value = 109
$ws1.Cells.Item(117, 5).Value = "LP1912"

$ws1.Cells.Item(118, 1).Value = "08:33:47"
$ws1.Cells.Item(118, 2).Value = "09:34"
$ws1.Cells.Item(118, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(118, 4).Value = 61
$ws1.Cells.Item(118, 5).Value = "LP1912"

$ws1.Cells.Item(119, 1).Value = "08:47:19"
$ws1.Cells.Item(119, 2).Value = "09:35"
$ws1.Cells.Item(119, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(119, 4).Value = 48
$ws1.Cells.Item(119, 5).Value = "LP1912"

$ws1.Cells.Item(120, 1).Value = "08:47:19"
$ws1.Cells.Item(120, 2).Value = "09:35"
$ws1.Cells.Item(120, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(120, 4).Value = 48
$ws1.Cells.Item(120, 5).Value = "LP1912"

$ws1.Cells.Item(121, 1).Value = "07:44:08"
$ws1.Cells.Item(121, 2).Value = "09:36"
$ws1.Cells.Item(121, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(121, 4).Value = 112
$ws1.Cells.Item(121, 5).Value = "LP1912"

$ws1.Cells.Item(122, 1).Value = "08:16:48"
$ws1.Cells.Item(122, 2).Value = "09:37"
$ws1.Cells.Item(122, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(122, 4).Value = 81
$ws1.Cells.Item(122, 5).Value = "LP1912"

$ws1.Cells.Item(123, 1).Value = "08:16:48"
$ws1.Cells.Item(123, 2).Value = "09:41"
$ws1.Cells.Item(123, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(123, 4).Value = 85
$ws1.Cells.Item(123, 5).Value = "LP1912"

$ws1.Cells.Item(124, 1).Value = "07:44:08"
$ws1.Cells.Item(124, 2).Value = "09:41"
$ws1.Cells.Item(124, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(124, 4).Value = 68
$ws1.Cells.Item(124, 5).Value = "LP1912"

$ws1.Cells.Item(125, 1).Value = "07:44:08"
$ws1.Cells.Item(125, 2).Value = "09:42"
$ws1.Cells.Item(125, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(125, 4).Value = 118
$ws1.Cells.Item(125, 5).Value = "LP1912"

$ws1.Cells.Item(126, 1).Value = "07:57:27"
$ws1.Cells.Item(126, 2).Value = "09:43"
$ws1.Cells.Item(126, 3).Value = "14_ABASTO"
$ws1.Cells.Item(126, 4).Value = 106
$ws1.Cells.Item(126, 5).Value = "LP1912"

$ws1.Cells.Item(127, 1).Value = "08:16:48"
$ws1.Cells.Item(127, 2).Value = "10:10"
$ws1.Cells.Item(127, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(127, 4).Value = 114
$ws1.Cells.Item(127, 5).Value = "LP1912"

$ws1.Cells.Item(128, 1).Value = "08:16:48"
$ws1.Cells.Item(128, 2).Value = "10:12"
$ws1.Cells.Item(128, 3).Value = "15_ABASTO"
$ws1.Cells.Item(128, 4).Value = 116
$ws1.Cells.Item(128, 5).Value = "LP1912"

$ws1.Cells.Item(129, 1).Value = "08:33:47"
$ws1.Cells.Item(129, 2).Value = "10:21"
$ws1.Cells.Item(129, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(129, 4).Value = 108
$ws1.Cells.Item(129, 5).Value = "LP1912"

$ws1.Cells.Item(130, 1).Value = "08:33:47"
$ws1.Cells.Item(130, 2).Value = "10:22"
$ws1.Cells.Item(130, 3).Value = "17_ROMERO"
$ws1.Cells.Item(130, 4).Value = 109
$ws1.Cells.Item(130, 5).Value = "LP1912"

$ws1.Cells.Item(131, 1).Value = "08:33:47"
$ws1.Cells.Item(131, 2).Value = "10:26"
$ws1.Cells.Item(131, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(131, 4).Value = 113
$ws1.Cells.Item(131, 5).Value = "LP1912"

$ws1.Cells.Item(132, 1).Value = "08:47:19"
$ws1.Cells.Item(132, 2).Value = "10:42"
$ws1.Cells.Item(132, 3).Value = "17_ROMERO"
$ws1.Cells.Item(132, 4).Value = 115
$ws1.Cells.Item(132, 5).Value = "LP1912"

$ws1.Cells.Item(133, 1).Value = "08:47:19"
$ws1.Cells.Item(133, 2).Value = "10:43"
$ws1.Cells.Item(133, 3).Value = "14_ABASTO"
$ws1.Cells.Item(133, 4).Value = 116
$ws1.Cells.Item(133, 5).Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 08:47:19"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 08:47:19"
$ws3.Range("A3").Value = "Total filas: 24"

$ws3.Cells.Item(27, 1).Value = "08:47:19"
$ws3.Cells.Item(27, 2).Value = "09:12"
$ws3.Cells.Item(27, 3).Value = "215D_LA PLATA"
$ws3.Cells.Item(27, 4).Value = 25
$ws3.Cells.Item(27, 5).Value = "L6203"

$ws3.Cells.Item(28, 1).Value = "08:16:48"
$ws3.Cells.Item(28, 2).Value = "10:02"
$ws3.Cells.Item(28, 3).Value = "215B_LP-P MOR-40 Y 115"
$ws3.Cells.Item(28, 4).Value = 106
$ws3.Cells.Item(28, 5).Value = "L6173"

$ws3.Cells.Item(29, 1).Value = "08:47:19"
$ws3.Cells.Item(29, 2).Value = "10:03"
$ws3.Cells.Item(29, 3).Value = "215B_LP-P MOR-40 Y 115"
$ws3.Cells.Item(29, 4).Value = 76
$ws3.Cells.Item(29, 5).Value = "L6173"

Write-Output "Edit complete"